# [DC2] adding player position column from 'Player Details' sheet to 'Prediction Stats' sheet
$wb = $excel.ActiveWorkbook
$wsDetails = $wb.Worksheets.Item("Player Details")
$wsPred = $wb.Worksheets.Item("Prediction Stats")

# Build a Name -> Position lookup table from the 'Player Details' sheet
# (column A = Name, column E = Position)
$lastRowDetails = $wsDetails.Cells.Item($wsDetails.Rows.Count, 1).End(-4162).Row
$positionByName = @{}
for ($r = 2; $r -le $lastRowDetails; $r++) {
    $playerName = $wsDetails.Cells.Item($r, 1).Value()
    $playerPosition = $wsDetails.Cells.Item($r, 5).Value()
    $positionByName[$playerName] = $playerPosition
}

# Populate column AD on 'Prediction Stats' with each player's position,
# looked up by matching the player's Name (column A) against 'Player Details'
$lastRowPred = $wsPred.Cells.Item($wsPred.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowPred; $r++) {
    $playerName = $wsPred.Cells.Item($r, 1).Value()
    $wsPred.Cells.Item($r, 30).Value = $positionByName[$playerName]
}
